$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "320.62"
Set-TextValue "E2" "-3.26%"
Set-TextValue "D3" "42.50"
Set-TextValue "E3" "-6.44%"
Set-TextValue "D4" "5.195"
Set-TextValue "E4" "-6.50%"
Set-TextValue "D5" "0.08173"
Set-TextValue "E5" "-2.18%"
Set-TextValue "D6" "4.314"
Set-TextValue "E6" "-3.04%"
Set-TextValue "D7" "1.813"
Set-TextValue "E7" "-13.42%"
Set-TextValue "D8" "0.9336"
Set-TextValue "E8" "-5.30%"
Set-TextValue "D9" "0.1111"
Set-TextValue "E9" "-7.59%"
Set-TextValue "D10" "0.1865"
Set-TextValue "E10" "-2.92%"
$ws.Range("B11").Value = "BitrueCoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue "D11" "0.04715"
Set-TextValue "E11" "1.06%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue "D12" "0.09406"
Set-TextValue "E12" "-4.58%"
Set-TextValue "D13" "7.421"
Set-TextValue "E13" "-28.03%"
Set-TextValue "E14" "-0.12%"
Set-TextValue "D15" "0.001309"
Set-TextValue "E15" "1.89%"
Set-TextValue "D16" "0.005791"
Set-TextValue "E16" "-1.85%"
$ws.Range("B17").Value = "HotbitToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb"
Set-TextValue "D17" "0.004298"
Set-TextValue "E17" "-5.16%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
Set-TextValue "D18" "3.356"
Set-TextValue "E18" "-1.07%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
Set-TextValue "D19" "2.533"
Set-TextValue "E19" "-0.52%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
Set-TextValue "D20" "0.3389"
Set-TextValue "E20" "1.49%"
$ws.Range("B21").Value = "ProBitToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
Set-TextValue "D21" "0.1391"
Set-TextValue "E21" "1.35%"
$ws.Range("B22").Value = "ZBToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
Set-TextValue "D22" "0.2549"
Set-TextValue "E22" "-0.61%"
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue "D23" "0.04149"
Set-TextValue "E23" "-0.10%"
$ws.Range("B24").Value = "BitKan"
$ws.Range("C24").Value = "https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan"
Set-TextValue "D24" "0.001245"
Set-TextValue "E24" "-3.77%"
Set-TextValue "E25" "-7.73%"
Set-TextValue "D26" "0.0002984"
Set-TextValue "E26" "-20.33%"
Set-TextValue "D38" "0.02702"
Set-TextValue "E38" "0.20%"
Set-TextValue "D39" "0.05547"
Set-TextValue "E39" "-3.43%"
Set-TextValue "D40" "0.008106"
Set-TextValue "E40" "2.76%"
Set-TextValue "D41" "0.1398"
Set-TextValue "E41" "-2.57%"
Set-TextValue "D42" "0.006558"
Set-TextValue "E42" "-13.23%"
Set-TextValue "D43" "0.002088"
Set-TextValue "E43" "-0.52%"
Set-TextValue "D44" "0.008268"
Set-TextValue "E44" "-7.48%"
Set-TextValue "D45" "0.3484"
Set-TextValue "E45" "2.32%"
Set-TextValue "D46" "0.00006934"
Set-TextValue "E46" "-1.51%"
Set-TextValue "E47" "-0.03%"
Set-TextValue "D48" "0.003359"
Set-TextValue "E48" "-2.69%"
Set-TextValue "D49" "0.003536"
Set-TextValue "E49" "0.01%"
Set-TextValue "E50" "-0.03%"
Set-TextValue "E51" "-0.03%"
